# Updated symbol list on Fri Dec 23 14:49:36 UTC 2022 with GitHub Actions
#
# The "Price" column (D) stores values as plain text (inline strings), even
# though they look numeric ("244.54", "5.388", ...). A plain
# `Range.Value = "..."` assignment would let Excel auto-convert those
# numeric-looking strings into real numbers (and normalise away trailing
# zeros, e.g. "5.390" -> 5.39). To keep them as text - exactly like the
# source data - we stage the new value in a scratch cell that has been
# explicitly formatted as Text ("@"), copy it, and paste it (values +
# formats) onto the destination cell. Because the destination cell's own
# number format is left as General (style untouched), the paste keeps the
# text as text without altering the destination cell's style index.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$helperAddr = "Z1000"

function Set-TextValue {
    param($ws, $addr, $val)

    $helper = $ws.Range($helperAddr)
    $helper.NumberFormat = "@"
    $helper.Value = $val
    $helper.Copy()
    $ws.Range($addr).PasteSpecial(-4104)  # xlPasteAll
    $helper.Clear()
}

# Plain text columns (Coin / Link / Volume label) - safe to assign directly,
# these values never look like numbers.
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("E41").Value = "40KickTokenKICK"

$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("E42").Value = "41BKEXTokenBKK"

$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("E43").Value = "42CEJICEJI"

$ws.Range("E18").Value = "17OneONEWorstin24h"

# Price column (D) - numeric-looking text, written via the text-preserving
# helper so it stays text (t="inlineStr"/shared-string) instead of becoming
# a Number cell.
Set-TextValue $ws "D2"  "244.75"
Set-TextValue $ws "D3"  "21.92"
Set-TextValue $ws "D4"  "5.390"
Set-TextValue $ws "D5"  "0.05850"
Set-TextValue $ws "D6"  "3.395"
Set-TextValue $ws "D7"  "6.356"
Set-TextValue $ws "D8"  "0.8162"
Set-TextValue $ws "D9"  "1.010"
Set-TextValue $ws "D10" "0.1424"
Set-TextValue $ws "D11" "0.03731"
Set-TextValue $ws "D12" "0.07497"
Set-TextValue $ws "D13" "0.03049"
Set-TextValue $ws "D14" "4.221"
Set-TextValue $ws "D15" "0.09389"
Set-TextValue $ws "D16" "0.001604"
Set-TextValue $ws "D17" "0.04835"
Set-TextValue $ws "D18" "0.0005899"
Set-TextValue $ws "D19" "0.006055"
Set-TextValue $ws "D21" "0.001001"
Set-TextValue $ws "D22" "0.0001501"
Set-TextValue $ws "D23" "3.692"
Set-TextValue $ws "D24" "2.220"
Set-TextValue $ws "D26" "0.1298"
Set-TextValue $ws "D27" "0.0002905"
Set-TextValue $ws "D40" "0.03852"
Set-TextValue $ws "D41" "0.006375"
Set-TextValue $ws "D42" "0.1072"
Set-TextValue $ws "D43" "0.002698"
Set-TextValue $ws "D44" "0.006243"
Set-TextValue $ws "D45" "0.00005631"
Set-TextValue $ws "D47" "0.8213"
Set-TextValue $ws "D48" "0.1422"
Set-TextValue $ws "D49" "0.00002102"
Set-TextValue $ws "D50" "0.01011"
